# Updated policies and graphs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new start-date (add_start) values for several policy rows.
$ws.Range("C13").Value = 43922
$ws.Range("C14").Value = 43922
$ws.Range("C19").Value = 43908
$ws.Range("C20").Value = 43908
$ws.Range("C21").Value = 43922
$ws.Range("C24").Value = 43908
$ws.Range("C27").Value = 43922
$ws.Range("C28").Value = 43922

# Apply the same date style (style index 2 in the original file) used by
# the sibling C-column cells (e.g. C3) so the new cells match formatting.
$ws.Range("C13").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C14").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C19").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C20").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C21").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C24").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C27").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("C28").NumberFormat = $ws.Range("C3").NumberFormat

# Select C24, matching the saved view/selection state.
$ws.Range("C24").Select()

# Autofit columns A:D so the widths reflect the newly entered data (Excel
# recomputed "best fit" widths for these columns when the file was last
# saved after the new dates were entered).
$ws.Columns("A:D").AutoFit() | Out-Null
